$d = $word.ActiveDocument

# --- 1) Grab the run layout of the "Meta description" paragraph (an empty
#        leading run followed by a bold run) - we reuse this exact run
#        structure for the new paragraph added near the end of the
#        document. NOTE: FormattedText is a live reference into the
#        document, so we must finish using it (paste it below) BEFORE we
#        delete its source paragraph, or it will "follow" whatever content
#        slides into that position after the delete. ---
$metaPara = $d.Paragraphs.Item(2)
$metaRunLayout = $metaPara.Range.FormattedText

# --- 2) Near the end of the document, insert a new bold paragraph with the
#        page title text, right before the final (italic) paragraph. ---
$count = $d.Paragraphs.Count
$priorPara = $d.Paragraphs.Item($count - 1)
$priorPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($count - 1 + 1)
$newPara.Style = "Normal"
$newPara.Range.FormattedText = $metaRunLayout

$pStart = $newPara.Range.Start
$pEnd = $newPara.Range.End
$fullRange = $d.Range($pStart, $pEnd - 1)
$fullRange.Text = "Play Fairy Dust Extreme Free - Review of Fantasy-themed Slot"
$fullRange.Font.Bold = 1

# --- 3) Now remove the "Meta description" paragraph that follows the title ---
$d.Paragraphs.Item(2).Range.Delete()

# --- 4) Replace the text of the last (italic) paragraph with the meta
#        description text, keeping its italic formatting intact. ---
$old = 'Create a feature image fitting for the game "Fairy Dust Extreme". The image should be in cartoon style and feature a happy Maya warrior wearing glasses. For the feature image, let''s incorporate both the fantasy fairy tale theme and the happy Maya warrior with glasses. We can have the Maya warrior surrounded by magical fairies, with a big smile on his face while holding a bag of lucky fairy dust. The fairies can have different colors to represent the game''s different bonus features, such as golden fairies for expanding Wilds, ruby fairies for medium paying symbols, and sapphire fairies for turning symbols into Wilds. The background can be a magical forest with mushrooms, plants, flowers, and small creatures to represent the game''s symbols. The overall design should be bright and vibrant, with a playful and enchanting tone to catch the attention of potential players.'
$new = 'Play Fairy Dust Extreme for free and read our review of this fantasy-themed online slot game. Find out the pros and cons of Fairy Dust Extreme.'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "Edit complete"
